# Better iterative FFT (#6) - updated timings
# Update the "Avg Time (ms)" column (D) on the Timings sheet with the
# new measured values, then leave the selection on D8 (matching the
# author's last-saved cursor position).
#
# Switch to manual calculation first so the PivotTable on the "Chart"
# sheet (still pointing at its last-refreshed cache) is not silently
# recomputed as a side effect of editing the source table - the author
# did not refresh the pivot/chart in this commit.
$excel.Calculation = -4135   # xlCalculationManual

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Timings")

# Recursive method timings (rows 2-8)
$ws.Range("D2").Value  = 0.00133333
$ws.Range("D3").Value  = 0.0303333
$ws.Range("D4").Value  = 0.586
$ws.Range("D5").Value  = 12.9337
$ws.Range("D6").Value  = 239.302
$ws.Range("D7").Value  = 4717.36
$ws.Range("D8").Value  = 128852

# Iterative method timings (rows 9-15)
$ws.Range("D9").Value  = 0.0004
$ws.Range("D10").Value = 0.0044
$ws.Range("D11").Value = 0.0995
$ws.Range("D12").Value = 1.9017
$ws.Range("D13").Value = 50.8885
$ws.Range("D14").Value = 1113.12
$ws.Range("D15").Value = 21817.8

# Match the saved cursor/selection position recorded in the workbook.
$ws.Activate()
$ws.Range("D8").Select()
